# Update scripts with new TPM values: recompute NATMI LR-pair stats for
# Rspo2-Rnf43 (Young D7) and add a new "Resolving-Mac" cluster to the
# sending/target cluster combinations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> MuSCs (values recomputed with new TPM)
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Rspo2"
$ws.Cells.Item(2, 3).Value = "Rnf43"
$ws.Cells.Item(2, 4).Value = "MuSCs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.597878666666666
$ws.Cells.Item(2, 8).Value = 4.793635999999999
$ws.Cells.Item(2, 9).Value = 0.8992131381376172
$ws.Cells.Item(2, 10).Value = 0.899213138137617
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.07557799999999999
$ws.Cells.Item(2, 14).Value = 0.226734
$ws.Cells.Item(2, 15).Value = 0.8607417137086825
$ws.Cells.Item(2, 16).Value = 0.8607417137086824
$ws.Cells.Item(2, 17).Value = 0.1207644738693333
$ws.Cells.Item(2, 18).Value = 1.086880264824
$ws.Cells.Item(2, 19).Value = 0.7739902575099348
$ws.Cells.Item(2, 20).Value = 0.7739902575099347

# Row 3: FAPs -> Resolving-Mac (new row, new target cluster)
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Rspo2"
$ws.Cells.Item(3, 3).Value = "Rnf43"
$ws.Cells.Item(3, 4).Value = "Resolving-Mac"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.597878666666666
$ws.Cells.Item(3, 8).Value = 4.793635999999999
$ws.Cells.Item(3, 9).Value = 0.8992131381376172
$ws.Cells.Item(3, 10).Value = 0.899213138137617
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.01222766666666667
$ws.Cells.Item(3, 14).Value = 0.036683
$ws.Cells.Item(3, 15).Value = 0.1392582862913176
$ws.Cells.Item(3, 16).Value = 0.1392582862913176
$ws.Cells.Item(3, 17).Value = 0.01953832770977778
$ws.Cells.Item(3, 18).Value = 0.175844949388
$ws.Cells.Item(3, 19).Value = 0.1252228806276824
$ws.Cells.Item(3, 20).Value = 0.1252228806276824

# Row 4: MuSCs -> MuSCs (new row)
$ws.Cells.Item(4, 1).Value = "MuSCs"
$ws.Cells.Item(4, 2).Value = "Rspo2"
$ws.Cells.Item(4, 3).Value = "Rnf43"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1790956666666667
$ws.Cells.Item(4, 8).Value = 0.537287
$ws.Cells.Item(4, 9).Value = 0.1007868618623829
$ws.Cells.Item(4, 10).Value = 0.1007868618623829
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.07557799999999999
$ws.Cells.Item(4, 14).Value = 0.226734
$ws.Cells.Item(4, 15).Value = 0.8607417137086825
$ws.Cells.Item(4, 16).Value = 0.8607417137086824
$ws.Cells.Item(4, 17).Value = 0.01353569229533333
$ws.Cells.Item(4, 18).Value = 0.121821230658
$ws.Cells.Item(4, 19).Value = 0.08675145619874775
$ws.Cells.Item(4, 20).Value = 0.08675145619874773

# Row 5: MuSCs -> Resolving-Mac (new row, new target cluster)
$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Rspo2"
$ws.Cells.Item(5, 3).Value = "Rnf43"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.1790956666666667
$ws.Cells.Item(5, 8).Value = 0.537287
$ws.Cells.Item(5, 9).Value = 0.1007868618623829
$ws.Cells.Item(5, 10).Value = 0.1007868618623829
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.01222766666666667
$ws.Cells.Item(5, 14).Value = 0.036683
$ws.Cells.Item(5, 15).Value = 0.1392582862913176
$ws.Cells.Item(5, 16).Value = 0.1392582862913176
$ws.Cells.Item(5, 17).Value = 0.002189922113444444
$ws.Cells.Item(5, 18).Value = 0.019709299021
$ws.Cells.Item(5, 19).Value = 0.0140354056636352
$ws.Cells.Item(5, 20).Value = 0.0140354056636352
